$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write one data row (A=index, B=name, C/D=bus refs, E=in_service bool)
# and make sure column A keeps the same visual style ("s=1": bold, centered,
# thin border) that every other row in column A already uses.
# NOTE: positional parameters only - named parameter binding is not
# supported by this interpreter.
function Set-DataRow($Row, $Index, $Name, $C, $D, $InService) {
    $a = $ws.Cells.Item($Row, 1)
    $a.Value = $Index
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108   # xlCenter
    $a.VerticalAlignment = -4160     # xlTop
    $a.Borders.LineStyle = 1         # xlContinuous (thin box border)

    $ws.Cells.Item($Row, 2).Value = $Name
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $InService
}

# Two new contingency lines ("line7", "line8") are inserted right after
# "line6", pushing the previously-existing "extr1".."extr8" rows down by two
# rows (old row 8 -> new row 10, ..., old row 15 -> new row 17). The index
# column (A) is renumbered sequentially, and several C/D/E values for the
# shifted rows change as part of the same commit.

Set-DataRow 8  6  "line7" 14 11 $false
Set-DataRow 9  7  "line8" 16 9  $true

Set-DataRow 10 8  "extr1" 5  12 $true
Set-DataRow 11 9  "extr2" 5  9  $true
Set-DataRow 12 10 "extr3" 10 11 $true
Set-DataRow 13 11 "extr4" 7  8  $false
Set-DataRow 14 12 "extr5" 9  11 $false
Set-DataRow 15 13 "extr6" 7  11 $false
Set-DataRow 16 14 "extr7" 5  7  $false
Set-DataRow 17 15 "extr8" 8  5  $true
